$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.945.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "'2.373.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'319.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.85%  "
$ws.Range("D6").Value = "'107.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'42.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "'16.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.97%  "
$ws.Range("D16").Value = "'2.733.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'2.394.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'42.928.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "'7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'76.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").Value = "'257.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.32%  "
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("D25").Value = "'9.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "'23.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").Value = "'171.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("D31").Value = "'36.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.59%  "
$ws.Range("E35").Value = "  +12.47%  "
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("D38").Value = "'0.0366"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'3.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "'72.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'12.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").Value = "'90.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").Value = "'113.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.71%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").Value = "'9.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "'77.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.17%  "
$ws.Range("D51").Value = "'1.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
